$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('T2').Value = 'maa://22742 (91.56), *maa://20791 (62.86)'
$ws.Range('H3').Value = 'maa://21247 (98.41), *maa://22748 (60.0)'
$ws.Range('L3').Value = '*maa://22880 (65.93), maa://20276 (84.81), *maa://22749 (72.73)'
$ws.Range('T4').Value = 'maa://32509 (97.96), maa://27295 (83.61), maa://22754 (91.67), *maa://21746 (55.81), *maa://31008 (78.05)'
$ws.Range('D5').Value = 'maa://21245 (83.18), maa://22744 (84.0)'
$ws.Range('L7').Value = 'maa://28624 (93.1), maa://24957 (97.62)'
$ws.Range('A8').Value = '更新日期：2024.12.19 13:19:06'
$ws.Range('AB9').Value = 'maa://28711 (87.0), ***maa://22740 (5.88), **maa://39938 (48.0), **maa://27377 (42.86), ***maa://25174 (19.05), maa://40166 (93.33)'
$ws.Range('T10').Value = 'maa://27395 (95.95), maa://22755 (87.5), **maa://22756 (40.91), ***maa://21737 (10.61)'
$ws.Range('T11').Value = 'maa://22747 (93.42), maa://22501 (98.57)'
$ws.Range('X12').Value = 'maa://22753 (91.36), *maa://21485 (76.87), maa://37962 (88.89)'
$ws.Range('AF12').Value = '*maa://28932 (78.63), *maa://20106 (63.96), *maa://22769 (64.29)'
$ws.Range('AF13').Value = '**maa://22737 (32.37), maa://39883 (90.74), *maa://39885 (56.0)'
$ws.Range('L14').Value = 'maa://26245 (96.43), maa://21288 (96.27), maa://39841 (95.0), maa://36682 (97.37)'
$ws.Range('T14').Value = 'maa://22521 (93.81), maa://42751 (100.0)'
$ws.Range('D15').Value = '*maa://22743 (76.96), maa://22734 (83.9), *maa://30808 (64.52), **maa://36048 (31.82)'
$ws.Range('P15').Value = 'maa://24762 (90.2), *maa://22727 (70.0)'
$ws.Range('AF15').Value = 'maa://21364 (80.84), *maa://22766 (70.91), *maa://36666 (78.57)'
$ws.Range('AF16').Value = '*maa://23911 (65.35), maa://27755 (92.5)'
$ws.Range('D17').Value = 'maa://21624 (82.86)'
$ws.Range('H17').Value = 'maa://22430 (88.65), maa://39599 (86.11)'
$ws.Range('T17').Value = '**maa://42324 (48.15)'
$ws.Range('D18').Value = 'maa://24570 (97.09)'
$ws.Range('L18').Value = 'maa://22466 (88.81), *maa://22732 (50.6)'
$ws.Range('AB19').Value = '*maa://30709 (63.05), *maa://36668 (55.84)'
$ws.Range('AB21').Value = 'maa://21443 (80.06), ***maa://23820 (29.82)'
$ws.Range('T22').Value = 'maa://38495 (84.62)'
$ws.Range('L23').Value = 'maa://39756 (93.98), maa://39875 (93.65)'
$ws.Range('P23').Value = 'maa://30587 (91.98), *maa://29748 (75.59), ***maa://29785 (16.42), *maa://37566 (71.43)'
$ws.Range('T23').Value = 'maa://24387 (81.58), maa://31212 (93.1)'
$ws.Range('X24').Value = 'maa://29988 (86.67), maa://23504 (93.26), **maa://22892 (39.58), *maa://25141 (77.6), maa://36663 (80.6), ***maa://22815 (23.08)'
$ws.Range('AF24').Value = 'maa://22523 (85.49), maa://36672 (80.77), maa://29910 (92.45), **maa://21440 (34.55)'
$ws.Range('T25').Value = 'maa://20109 (92.4), maa://22545 (100.0), maa://42915 (100.0)'
$ws.Range('X26').Value = 'maa://24389 (96.3)'
$ws.Range('X28').Value = 'maa://39929 (89.68), ***maa://39723 (14.29), maa://41749 (90.57)'
$ws.Range('L31').Value = 'maa://35926 (93.96), maa://36258 (83.16), *maa://43904 (71.43)'
$ws.Range('H32').Value = 'maa://21895 (97.25), maa://36667 (98.33), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range('T32').Value = 'maa://42859 (96.1), maa://41108 (87.76), maa://41238 (96.15)'
$ws.Range('L35').Value = 'maa://41296 (96.55)'
$ws.Range('P37').Value = 'maa://21280 (89.0), *maa://21239 (66.67)'
$ws.Range('P38').Value = '*maa://24383 (67.71)'
$ws.Range('T38').Value = 'maa://30713 (96.77)'
$ws.Range('P39').Value = 'maa://24709 (92.13)'
$ws.Range('P40').Value = 'maa://23278 (95.96), maa://21386 (95.74), maa://36664 (90.74)'
$ws.Range('H43').Value = 'maa://22525 (92.65), maa://21284 (83.72)'
$ws.Range('H45').Value = 'maa://21229 (84.78), maa://30807 (95.45), *maa://22767 (55.0), ***maa://20796 (13.79), *maa://42459 (71.43)'
$ws.Range('P45').Value = '*maa://36237 (64.29)'
$ws.Range('T45').Value = '**maa://39364 (35.0)'
$ws.Range('H47').Value = 'maa://27410 (96.18), maa://29661 (97.83), maa://28038 (84.62)'
$ws.Range('H58').Value = '*maa://37964 (55.17)'
$ws.Range('H60').Value = '*maa://40438 (59.09)'
